$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first row of the sheet (the header row containing "Tanggal" and
# the day numbers 1-31). This shifts every remaining row up by one, which is
# why the sheet's used range shrinks from A1:AF205 to A1:AF204 and the
# "Tanggal" shared string ends up unused (and gets pruned on save).
$ws.Rows("1").Delete()

# Reflect where the user's selection ended up after the edit.
[void]$ws.Range("L7").Select()
